# "Group Reflection OneDrive ReSync"
#
# The captured OOXML diff for this commit adds no visible body text —
# it only adds the customXml parts (content-type schema, SharePoint
# form templates, and document-management properties) that OneDrive /
# SharePoint stamps onto a document the first time it resyncs a
# library item. Reproduce that via the Word object model by adding the
# same CustomXMLParts to the document.

$d = $word.ActiveDocument

# customXml/item1.xml - SharePoint content-type schema (columns,
# core-properties mapping, InfoPath "Person" field schema, etc.)
$contentTypeSchemaXml = @'
<?xml version="1.0" encoding="utf-8"?>
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x01010027B292595FC22049A95F0029909BF2E4" ma:contentTypeVersion="9" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="c077942260d98140fd3c47e4fedf63bd">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="e9d81c5b-9174-487f-9c97-32e1760a833a" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="46ac66f369aaff067c698ca5bb15b602" ns2:_="">
    <xsd:import namespace="e9d81c5b-9174-487f-9c97-32e1760a833a"/>
    <xsd:element name="properties">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element name="documentManagement">
            <xsd:complexType>
              <xsd:all>
                <xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceAutoKeyPoints" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceKeyPoints" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceAutoTags" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceOCR" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceGenerationTime" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceDateTaken" minOccurs="0"/>
              </xsd:all>
            </xsd:complexType>
          </xsd:element>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="e9d81c5b-9174-487f-9c97-32e1760a833a" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceAutoKeyPoints" ma:index="10" nillable="true" ma:displayName="MediaServiceAutoKeyPoints" ma:hidden="true" ma:internalName="MediaServiceAutoKeyPoints" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceKeyPoints" ma:index="11" nillable="true" ma:displayName="KeyPoints" ma:internalName="MediaServiceKeyPoints" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceAutoTags" ma:index="12" nillable="true" ma:displayName="Tags" ma:internalName="MediaServiceAutoTags" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceOCR" ma:index="13" nillable="true" ma:displayName="Extracted Text" ma:internalName="MediaServiceOCR" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceGenerationTime" ma:index="14" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceEventHashCode" ma:index="15" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceDateTaken" ma:index="16" nillable="true" ma:displayName="MediaServiceDateTaken" ma:hidden="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd" targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all">
    <xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/>
    <xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/>
    <xsd:element name="coreProperties" type="CT_coreProperties"/>
    <xsd:complexType name="CT_coreProperties">
      <xsd:all>
        <xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/>
        <xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/>
        <xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string">
          <xsd:annotation>
            <xsd:documentation>
                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.
                    </xsd:documentation>
          </xsd:annotation>
        </xsd:element>
        <xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/>
      </xsd:all>
    </xsd:complexType>
  </xsd:schema>
  <xs:schema xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema" targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified">
    <xs:element name="Person">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:DisplayName" minOccurs="0"/>
          <xs:element ref="pc:AccountId" minOccurs="0"/>
          <xs:element ref="pc:AccountType" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="DisplayName" type="xs:string"/>
    <xs:element name="AccountId" type="xs:string"/>
    <xs:element name="AccountType" type="xs:string"/>
    <xs:element name="BDCAssociatedEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
        <xs:attribute ref="pc:EntityNamespace"/>
        <xs:attribute ref="pc:EntityName"/>
        <xs:attribute ref="pc:SystemInstanceName"/>
        <xs:attribute ref="pc:AssociationName"/>
      </xs:complexType>
    </xs:element>
    <xs:attribute name="EntityNamespace" type="xs:string"/>
    <xs:attribute name="EntityName" type="xs:string"/>
    <xs:attribute name="SystemInstanceName" type="xs:string"/>
    <xs:attribute name="AssociationName" type="xs:string"/>
    <xs:element name="BDCEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:EntityDisplayName" minOccurs="0"/>
          <xs:element ref="pc:EntityInstanceReference" minOccurs="0"/>
          <xs:element ref="pc:EntityId1" minOccurs="0"/>
          <xs:element ref="pc:EntityId2" minOccurs="0"/>
          <xs:element ref="pc:EntityId3" minOccurs="0"/>
          <xs:element ref="pc:EntityId4" minOccurs="0"/>
          <xs:element ref="pc:EntityId5" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="EntityDisplayName" type="xs:string"/>
    <xs:element name="EntityInstanceReference" type="xs:string"/>
    <xs:element name="EntityId1" type="xs:string"/>
    <xs:element name="EntityId2" type="xs:string"/>
    <xs:element name="EntityId3" type="xs:string"/>
    <xs:element name="EntityId4" type="xs:string"/>
    <xs:element name="EntityId5" type="xs:string"/>
    <xs:element name="Terms">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermInfo">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermName" minOccurs="0"/>
          <xs:element ref="pc:TermId" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermName" type="xs:string"/>
    <xs:element name="TermId" type="xs:string"/>
  </xs:schema>
</ct:contentTypeSchema>
'@

# customXml/item2.xml - document-library form template pointers
$formTemplatesXml = @'
<?xml version="1.0" encoding="utf-8"?>
<?mso-contentType ?>
<FormTemplates xmlns="http://schemas.microsoft.com/sharepoint/v3/contenttype/forms">
  <Display>DocumentLibraryForm</Display>
  <Edit>DocumentLibraryForm</Edit>
  <New>DocumentLibraryForm</New>
</FormTemplates>
'@

# customXml/item3.xml - empty documentManagement properties bag
$propertiesXml = @'
<?xml version="1.0" encoding="utf-8"?>
<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
  <documentManagement/>
</p:properties>
'@

$null = $d.CustomXMLParts.Add($contentTypeSchemaXml)
$null = $d.CustomXMLParts.Add($formTemplatesXml)
$null = $d.CustomXMLParts.Add($propertiesXml)
